$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the X/Y coordinate columns
$ws.Range("J1").Value = "X_COORD"
$ws.Range("K1").Value = "Y_COORD"

# New coordinate data rows
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 500

$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 600

$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 700

$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 800

$ws.Range("J6").Value = 700
$ws.Range("K6").Value = 900

# Store the new coordinates as text-formatted, left-aligned numbers
# (matches the new style added to cellXfs: numFmtId 49, left alignment)
$ws.Range("J2:K6").HorizontalAlignment = -4131
$ws.Range("J2:K6").NumberFormat = "@"

# Column width adjustments (existing column C narrowed, new columns E:K sized)
$ws.Range("C1").ColumnWidth = 29.666666666666668
$ws.Range("E1").ColumnWidth = 19.666666666666668
$ws.Range("F1").ColumnWidth = 14.333333333333332
$ws.Range("G1").ColumnWidth = 32.166666666666664
$ws.Range("H1").ColumnWidth = 24.666666666666668
$ws.Range("I1").ColumnWidth = 52.666666666666664
$ws.Range("J1").ColumnWidth = 17.833333333333336
$ws.Range("K1").ColumnWidth = 12.666666666666666

# Selection now covers the freshly entered coordinate block
$null = $ws.Range("J2:K6").Select()
